# Composer CMS v2.0.beta6 (2015-04-02)  ->  v2.0.beta7 (2015-05-10)
# Updates the Date-styled paragraph's version/date text.

$d = $word.ActiveDocument

$oldVersion = "v2.0.beta6"
$newVersion = "v2.0.beta7"
$oldDate    = "(2015-04-02)"
$newDate    = "(2015-05-10)"

# --- locate the three runs (version / space / date) in the Date paragraph ---
$verRng = $d.Content
$verRng.Find.Execute($oldVersion, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verStart = $verRng.Start
$verEnd   = $verRng.End

$dateRng = $d.Content
$dateRng.Find.Execute($oldDate, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateStart = $dateRng.Start
$dateEnd   = $dateRng.End

# --- replace the whole "version space date" span with the new text in one go ---
$fullRng = $d.Range($verStart, $dateEnd)
$fullRng.Text = $newVersion + " " + $newDate

# --- recompute the boundary between the new version text and the trailing text ---
$newVerEnd    = $verStart + $newVersion.Length
$newDateStart = $newVerEnd + 1
$newDateEnd   = $newDateStart + $newDate.Length

# --- re-split the merged run back into version / space / date runs so the
#     paragraph keeps its original three-run shape (Word merges identically
#     formatted adjacent runs whenever a Range's .Text is reassigned; toggling
#     a character property and reverting it forces Word to re-break the run
#     at that boundary without altering the visible formatting). ---
$splitA = $d.Range($newVerEnd, $newDateStart)
$splitA.Bold = 1
$splitA.Bold = 0

$splitB = $d.Range($newDateStart, $newDateEnd)
$splitB.Bold = 1
$splitB.Bold = 0

Write-Output "Updated Composer CMS version/date line."
